# Update quizvragen via Admin
#
# 1) "DC" sheet: a new quiz question row was appended (row 4).
# 2) "Wiskunde 3" sheet: the erroneous duplicate first question row
#    (row 2, which still held leftover "stroom I" / Ohm's-law content)
#    was removed, shifting every following question up by one row.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new question to the "DC" sheet ---------------------
$dc = $wb.Worksheets.Item("DC")

$dc.Range("B4").Value = "mc"
$dc.Range("D4").Value = "Is dit een goede nieuwe vraag??"
$dc.Range("E4").Value = "['A. test 1', ' B. Test 2', ' C. Test 3']"
$dc.Range("F4").Value = 1

# --- 2. Remove the stale duplicate row from "Wiskunde 3" ---------------
$wisk = $wb.Worksheets.Item("Wiskunde 3")
$wisk.Range("A2").EntireRow.Delete()
